# trafo_id -> gridnode_id refactor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell J1 from "trafo_id" to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Update the selected cell in the sheet view (was N12, now G6)
$ws.Range("G6").Select()
